# "Ultimas pruebas realizadas con IACK. Modificacion codigo IACK2"
# Adds a new worksheet "Retardos4" (after the existing "Retardos3") holding
# a fourth batch of latency measurements, following the exact same layout
# as "Retardos3" (three side-by-side tables: t1 / t2(uS) / t3(us), with a
# bold header per block and an AVERAGE() row at the bottom).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new sheet at the end of the workbook and name it.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet  = $wb.Worksheets.Item($sheetCount)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Retardos4"

# ---------------------------------------------------------------------
# 2. Block titles (row 5).
# ---------------------------------------------------------------------
$ws4.Range("A5").Value = "Payload 12 bytes delay 10ms sin procesamiento"
$ws4.Range("E5").Value = "Payload 12 bytes delay 4ms con procesamiento"
$ws4.Range("I5").Value = "Payload 12 bytes delay 4ms sin procesamiento"

# ---------------------------------------------------------------------
# 3. Column headers (row 7).
# ---------------------------------------------------------------------
$ws4.Range("A7").Value = "t1"
$ws4.Range("B7").Value = "t2(uS)"
$ws4.Range("C7").Value = "t3(us)"
$ws4.Range("E7").Value = "t1"
$ws4.Range("F7").Value = "t2(uS)"
$ws4.Range("G7").Value = "t3(us)"
$ws4.Range("I7").Value = "t1"
$ws4.Range("J7").Value = "t2(uS)"
$ws4.Range("K7").Value = "t3(us)"

# ---------------------------------------------------------------------
# 4. Data rows 8-17.
# ---------------------------------------------------------------------
$ws4.Range("A8").Value = 0
$ws4.Range("B8").Value = 2787
$ws4.Range("C8").Value = 16083
$ws4.Range("E8").Value = 0
$ws4.Range("F8").Value = 2796
$ws4.Range("G8").Value = 6092
$ws4.Range("I8").Value = 0
$ws4.Range("J8").Value = 3425
$ws4.Range("K8").Value = 4913

$ws4.Range("A9").Value = 0
$ws4.Range("B9").Value = 3750
$ws4.Range("C9").Value = 16722
$ws4.Range("E9").Value = 0
$ws4.Range("F9").Value = 1835
$ws4.Range("G9").Value = 8014
$ws4.Range("I9").Value = 0
$ws4.Range("J9").Value = 2478
$ws4.Range("K9").Value = 7788

$ws4.Range("A10").Value = 0
$ws4.Range("B10").Value = 1525
$ws4.Range("C10").Value = 16085
$ws4.Range("E10").Value = 0
$ws4.Range("F10").Value = 3748
$ws4.Range("G10").Value = 7373
$ws4.Range("I10").Value = 0
$ws4.Range("J10").Value = 2466
$ws4.Range("K10").Value = 6189

$ws4.Range("A11").Value = 0
$ws4.Range("B11").Value = 2151
$ws4.Range("C11").Value = 18324
$ws4.Range("E11").Value = 0
$ws4.Range("F11").Value = 3122
$ws4.Range("G11").Value = 6730
$ws4.Range("I11").Value = 0
$ws4.Range("J11").Value = 1524
$ws4.Range("K11").Value = 6508

$ws4.Range("A12").Value = 0
$ws4.Range("B12").Value = 2796
$ws4.Range("C12").Value = 17686
$ws4.Range("E12").Value = 0
$ws4.Range("F12").Value = 3436
$ws4.Range("G12").Value = 7691
$ws4.Range("I12").Value = 0
$ws4.Range("J12").Value = 2153
$ws4.Range("K12").Value = 7780

$ws4.Range("A13").Value = 0
$ws4.Range("B13").Value = 2153
$ws4.Range("C13").Value = 17363
$ws4.Range("E13").Value = 0
$ws4.Range("F13").Value = 2146
$ws4.Range("G13").Value = 5779
$ws4.Range("I13").Value = 0
$ws4.Range("J13").Value = 1518
$ws4.Range("K13").Value = 5556

$ws4.Range("A14").Value = 0
$ws4.Range("B14").Value = 1520
$ws4.Range("C14").Value = 16087
$ws4.Range("E14").Value = 0
$ws4.Range("F14").Value = 3748
$ws4.Range("G14").Value = 7375
$ws4.Range("I14").Value = 0
$ws4.Range("J14").Value = 1830
$ws4.Range("K14").Value = 5876

$ws4.Range("A15").Value = 0
$ws4.Range("B15").Value = 3745
$ws4.Range("C15").Value = 16080
$ws4.Range("E15").Value = 0
$ws4.Range("F15").Value = 2476
$ws4.Range("G15").Value = 5451
$ws4.Range("I15").Value = 0
$ws4.Range("J15").Value = 1514
$ws4.Range("K15").Value = 6190

$ws4.Range("A16").Value = 0
$ws4.Range("B16").Value = 2167
$ws4.Range("C16").Value = 16091
$ws4.Range("E16").Value = 0
$ws4.Range("F16").Value = 2467
$ws4.Range("G16").Value = 6730
$ws4.Range("I16").Value = 0
$ws4.Range("J16").Value = 3434
$ws4.Range("K16").Value = 6189

$ws4.Range("A17").Value = 0
$ws4.Range("B17").Value = 3114
$ws4.Range("C17").Value = 16399
$ws4.Range("E17").Value = 0
$ws4.Range("F17").Value = 1514
$ws4.Range("G17").Value = 5460
$ws4.Range("I17").Value = 0
$ws4.Range("J17").Value = 1510
$ws4.Range("K17").Value = 6499

# ---------------------------------------------------------------------
# 5. Averages row (18).
# ---------------------------------------------------------------------
$ws4.Range("A18").Value = "Promedio"
$ws4.Range("B18").Formula = "=AVERAGE(B8:B17)"
$ws4.Range("C18").Formula = "=AVERAGE(C8:C17)"
$ws4.Range("E18").Value = "Promedio"
$ws4.Range("F18").Formula = "=AVERAGE(F8:F17)"
$ws4.Range("G18").Formula = "=AVERAGE(G8:G17)"
$ws4.Range("I18").Value = "Promedio"
$ws4.Range("J18").Formula = "=AVERAGE(J8:J17)"
$ws4.Range("K18").Formula = "=AVERAGE(K8:K17)"

# ---------------------------------------------------------------------
# 6. Bold the title/header/label cells (matches the style used on the
#    other "Retardos" sheets).
# ---------------------------------------------------------------------
$boldRanges = "A5", "E5", "I5", "A7:C7", "E7:G7", "I7:K7", "A18", "E18", "I18"
foreach ($addr in $boldRanges) {
    $ws4.Range($addr).Font.Bold = $true
}

# ---------------------------------------------------------------------
# 7. View state: Retardos3 is no longer the active tab - restore its
#    plain selection, then make Retardos4 the active/visible sheet with
#    its own selection, matching the saved workbook state.
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Retardos3")
[void]$ws3.Range("A5:C18").Select()

[void]$ws4.Activate()
[void]$ws4.Range("I18").Select()
